$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "green running leggings for women",
    "green running skirt",
    "green running tights",
    "green running tights men",
    "green sheer tights",
    "green shiny tights",
    "green thick tights",
    "green tights for women",
    "green tights for women 3x",
    "green vibe weight loss",
    "green zebra leggings",
    "grey addias leggings",
    "grey adidas leggings women",
    "grey athletic tape",
    "grey cigarette pants",
    "grey compression leggings women",
    "grey compression pants",
    "grey stability ball",
    "grey thigh high tights",
    "grip tape plan b",
    "grooming length chart",
    "gt basketball",
    "guess the basketball team",
    "guide pro pants",
    "gym apparel for women",
    "gym capri",
    "gym capri leggings women high waist",
    "gym capris women",
    "gym close for women",
    "gym clothes for women shorts",
    "gym clothing women",
    "gym compression leggings women",
    "gym compression shorts",
    "gym floor tape",
    "gym knee support",
    "gym legging",
    "gym legging for women",
    "gym leggings",
    "gym leggings for women",
    "gym leggings with pocket",
    "gym leggings women",
    "gym lower",
    "gym lowers for women",
    "gym pant",
    "gym pants",
    "gym pants for women",
    "gym pants for women pack",
    "gym pants women",
    "gym pants women capri",
    "gym people leggings for women",
    "gym shorts yoga",
    "gym tape floor",
    "gym tights",
    "gym tights for men",
    "gym tights for women",
    "gym tights for women with pockets",
    "gym tights with pockets",
    "gym tights women",
    "gym trousers for women",
    "gym underwear women",
    "gym yoga leggings",
    "gym yoga pants",
    "gym yoga pants for women",
    "gym yoga women",
    "gymshark bottoms",
    "haha ball track",
    "half a cow",
    "half green half red pants",
    "half leggings for women",
    "half marathon gear",
    "half marathon tank",
    "half pants",
    "half pants for ladies",
    "half pants for women",
    "half pants for women for summer",
    "half running tights men",
    "half skirt",
    "half tights",
    "half tights running",
    "half tights women",
    "half trousers",
    "hamstring ball",
    "hamstring brace",
    "hamstring brace compression",
    "hamstring brace for men",
    "hamstring brace for women",
    "hamstring brace women",
    "hamstring braces",
    "hamstring compression",
    "hamstring compression brace",
    "hamstring compression men",
    "hamstring compression pants",
    "hamstring compression shorts",
    "hamstring compression shorts for men",
    "hamstring compression shorts men",
    "hamstring compression shorts women",
    "hamstring compression women",
    "hamstring equipment",
    "hamstring exercise",
    "hamstring exercise equipment"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
